$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text-formatted numbers (e.g. "60.780.87", subscript-notation
# "0.0[3]0868", trailing zeros that must be preserved, "." used as a thousands
# separator, etc.). Force Text format before assigning so Excel does not
# auto-convert them to native numbers (which would strip formatting/precision),
# then restore General so the stored cell style matches the source workbook.
$dCells = @("D2", "D3", "D5", "D6", "D9", "D11", "D12", "D13", "D15", "D16", "D17", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D33", "D34", "D35", "D36", "D37", "D41", "D42", "D43", "D44", "D45", "D47", "D50", "D51")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "60.780.87"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "2.906.42"
$ws.Range("E3").Value = "  -2.84%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "526.55"
$ws.Range("E5").Value = "  -2.93%  "
$ws.Range("D6").Value = "144.46"
$ws.Range("E6").Value = "  -4.89%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -4.07%  "
$ws.Range("D9").Value = "2.915.85"
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("E10").Value = "  -5.27%  "
$ws.Range("D11").Value = "6.14"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "0.358"
$ws.Range("E12").Value = "  -3.07%  "
$ws.Range("D13").Value = "3.413.82"
$ws.Range("E13").Value = "  -2.78%  "
$ws.Range("E14").Value = "  +2.52%  "
$ws.Range("D15").Value = "60.808.07"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("D16").Value = "22.55"
$ws.Range("E16").Value = "  -6.15%  "
$ws.Range("D17").Value = "2.916.77"
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("E18").Value = "  -3.87%  "
$ws.Range("E19").Value = "  -5.22%  "
$ws.Range("E20").Value = "  -4.37%  "
$ws.Range("D21").Value = "353.70"
$ws.Range("E21").Value = "  -6.66%  "
$ws.Range("D22").Value = "6.51"
$ws.Range("E22").Value = "  -3.16%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "5.70"
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("D25").Value = "64.81"
$ws.Range("E25").Value = "  -2.01%  "
$ws.Range("D26").Value = "0.450"
$ws.Range("E26").Value = "  -4.12%  "
$ws.Range("D27").Value = "0.179"
$ws.Range("E27").Value = "  -4.59%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0868"
$ws.Range("E29").Value = "  -7.10%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "7.81"
$ws.Range("E30").Value = "  -5.18%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  -2.72%  "
$ws.Range("D33").Value = "19.65"
$ws.Range("E33").Value = "  -4.05%  "
$ws.Range("D34").Value = "153.12"
$ws.Range("E34").Value = "  -4.98%  "
$ws.Range("D35").Value = "4.40"
$ws.Range("E35").Value = "  -4.27%  "
$ws.Range("D36").Value = "5.59"
$ws.Range("E36").Value = "  -6.71%  "
$ws.Range("D37").Value = "0.997"
$ws.Range("E37").Value = "  -6.76%  "
$ws.Range("E38").Value = "  -6.42%  "
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("E40").Value = "  -5.52%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.293.10"
$ws.Range("E41").Value = "  -5.21%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "0.652"
$ws.Range("E42").Value = "  -3.05%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "3.69"
$ws.Range("E43").Value = "  -5.50%  "
$ws.Range("D44").Value = "0.0583"
$ws.Range("E44").Value = "  -1.66%  "
$ws.Range("D45").Value = "20.36"
$ws.Range("E45").Value = "  -8.00%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "4.94"
$ws.Range("E47").Value = "  -4.50%  "
$ws.Range("E48").Value = "  -3.40%  "
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("D50").Value = "0.0914"
$ws.Range("E50").Value = "  -4.19%  "
$ws.Range("D51").Value = "18.55"
$ws.Range("E51").Value = "  -6.10%  "

foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "General" }
